$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.746.96'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '2.903.16'
$ws.Range('E3').Value = '  -2.82%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.54'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.31'
$ws.Range('E6').Value = '  -5.99%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.543'
$ws.Range('E8').Value = '  -4.99%  '
$ws.Range('D9').Value = '2.910.83'
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -5.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.05'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').Value = '  -3.46%  '
$ws.Range('D13').Value = '3.413.12'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('D15').Value = '60.717.40'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.53'
$ws.Range('E16').Value = '  -5.75%  '
$ws.Range('D17').Value = '2.914.48'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('E18').Value = '  -4.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.88'
$ws.Range('E19').Value = '  -5.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.58'
$ws.Range('E20').Value = '  -3.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.01'
$ws.Range('E21').Value = '  -7.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.47'
$ws.Range('E22').Value = '  -3.84%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.75'
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.85'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.450'
$ws.Range('E26').Value = '  -4.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.177'
$ws.Range('E27').Value = '  -6.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -4.85%  '
$ws.Range('D30').Value = '0.0₃0852'
$ws.Range('E30').Value = '  -9.81%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  -2.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.59'
$ws.Range('E33').Value = '  -4.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.03'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('E35').Value = '  -4.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.55'
$ws.Range('E36').Value = '  -6.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  -7.50%  '
$ws.Range('E38').Value = '  -6.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.56'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -5.38%  '
$ws.Range('E41').Value = '  -5.36%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.651'
$ws.Range('E42').Value = '  -3.58%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.283.44'
$ws.Range('E43').Value = '  -5.57%  '
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.27'
$ws.Range('E45').Value = '  -8.13%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('E47').Value = '  -4.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0236'
$ws.Range('E48').Value = '  -3.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.33'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0914'
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.34'
$ws.Range('E51').Value = '  -7.78%  '
